$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.872.95'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.622.73'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.34%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.42'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -0.08%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('E7').Value = '  +0.29%  '
$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.91'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  -2.16%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.257'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  +0.31%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0603'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.854.46'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '1.629.11'
$ws.Range('E13').Value = '  +0.19%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.95'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  -2.05%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.551'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  -2.00%  '
$__style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.47'
$ws.Range('D16').Style = $__style
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '27.866.96'
$ws.Range('E17').Value = '  -0.20%  '
$__style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '226.75'
$ws.Range('D18').Style = $__style
$ws.Range('E18').Value = '  -1.74%  '
$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.57'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').Value = '0.0₃0713'
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('E21').Value = '  +0.31%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.32'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -0.64%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.90'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('E24').Value = '  +1.35%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.12'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -0.49%  '
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.89'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.110'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  -1.40%  '
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.33'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  -0.44%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0479'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  -0.63%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.38'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('D33').Value = '1.419.28'
$ws.Range('E33').Value = '  +1.21%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.06'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -0.64%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +1.22%  '
$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.974'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('E39').Value = '  -0.70%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.844'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('E42').Value = '  -2.04%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.95'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  -1.79%  '
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.36'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  -2.79%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.77'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').Value = '1.763.01'
$ws.Range('E46').Value = '  -0.37%  '
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  -3.11%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.20'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  -1.91%  '
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0996'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('E51').Value = '  -0.41%  '
